$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Q2").Formula = '=CONCAT(P2,"_",B2,"_",D2,"_",E2)'
$f = $ws.Range("Q2").Formula
Write-Host "Q2 formula:" $f
$v = $ws.Range("Q2").Value()
Write-Host "Q2 value:" $v
